$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("A10").Value = "B0BF67DM6K"
$ws.Range("B10").Value = "SteelSeries Apex Pro TKL Wireless HyperMagnetic Gaming Keyboard — World's Fastest Keyboard — Esports Tenkeyless — OLED Screen — Adjustable Actuation — PBT Keycaps — Bluetooth — 2.4GHz — USB-C"
$ws.Range("C10").Value = 993.49

# Keep the date as literal text (matches "02/03/2024" inline string in the
# source row), preventing Excel from auto-converting it to a date serial.
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "02/03/2024"
$ws.Range("D10").Style = "Normal"
